$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.609.04"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "2.341.61"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.11"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.11"
$ws.Range("E6").Value = "  -2.11%  "
$ws.Range("E7").Value = "  -2.05%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.35"
$ws.Range("E10").Value = "  -2.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "51.80"
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0797"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("E14").Value = "  -2.44%  "
$ws.Range("D15").Value = "2.706.62"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.59"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "2.307.38"
$ws.Range("E17").Value = "  -2.73%  "
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "43.523.58"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.81"
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("D21").Value = "0.0₃0909"
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.13"
$ws.Range("E22").Value = "  -2.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.12"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "238.73"
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("E25").Value = "  -3.32%  "
$ws.Range("E26").Value = "  -2.94%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.09"
$ws.Range("E28").Value = "  -2.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("E29").Value = "  +4.07%  "
$ws.Range("E30").Value = "  -5.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.28"
$ws.Range("E31").Value = "  -3.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.45"
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  -3.90%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.42"
$ws.Range("E35").Value = "  -4.85%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.54"
$ws.Range("E36").Value = "  -4.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.09"
$ws.Range("E37").Value = "  -6.51%  "
$ws.Range("E38").Value = "  -4.34%  "
$ws.Range("E39").Value = "  -7.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.83"
$ws.Range("E40").Value = "  -5.32%  "
$ws.Range("E41").Value = "  -2.59%  "
$ws.Range("E42").Value = "  -2.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.42"
$ws.Range("E43").Value = "  -2.02%  "
$ws.Range("D44").Value = "1.990.34"
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("E45").Value = "  -1.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.61"
$ws.Range("E46").Value = "  -6.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.95"
$ws.Range("E47").Value = "  -6.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.95"
$ws.Range("E48").Value = "  -4.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "56.39"
$ws.Range("E49").Value = "  -2.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.92"
$ws.Range("E50").Value = "  +4.55%  "
$ws.Range("D51").Value = "2.568.36"
$ws.Range("E51").Value = "  +0.58%  "
